$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, capture the existing (soon to be modified) row 47 values so we can
# duplicate them into the new row 48 exactly as they were before the edit.
$a47 = $ws.Cells.Item(47,1).Value2
$b47 = $ws.Cells.Item(47,2).Value2
$c47 = $ws.Cells.Item(47,3).Value2
$d47 = $ws.Cells.Item(47,4).Value2
$e47 = $ws.Cells.Item(47,5).Value2
$f47 = $ws.Cells.Item(47,6).Value2
$g47 = $ws.Cells.Item(47,7).Value2
$h47 = $ws.Cells.Item(47,8).Value2
$i47 = $ws.Cells.Item(47,9).Value2
$j47 = $ws.Cells.Item(47,10).Value2
$k47 = $ws.Cells.Item(47,11).Value2
$l47 = $ws.Cells.Item(47,12).Value2
$m47 = $ws.Cells.Item(47,13).Value2
$n47 = $ws.Cells.Item(47,14).Value2
$o47 = $ws.Cells.Item(47,15).Value2
$p47 = $ws.Cells.Item(47,16).Value2
$q47 = $ws.Cells.Item(47,17).Value2
$r47 = $ws.Cells.Item(47,18).Value2

# Update row 47: new date (45135 -> 45265) and new volume (70 -> 100)
$ws.Cells.Item(47,4).Value = 45265
$ws.Cells.Item(47,10).Value = 100

# Add a new row 48, matching the previous (pre-edit) content of row 47.
$ws.Cells.Item(48,1).Value = $a47
$ws.Cells.Item(48,2).Value = $b47
$ws.Cells.Item(48,3).Value = $c47

$ws.Cells.Item(48,4).NumberFormat = $ws.Cells.Item(47,4).NumberFormat
$ws.Cells.Item(48,4).Value = $d47

$ws.Cells.Item(48,5).Value = $e47
$ws.Cells.Item(48,6).Value = $f47
$ws.Cells.Item(48,7).Value = $g47
$ws.Cells.Item(48,8).Value = $h47
$ws.Cells.Item(48,9).Value = $i47
$ws.Cells.Item(48,10).Value = $j47
$ws.Cells.Item(48,11).Value = $k47
$ws.Cells.Item(48,12).Value = $l47
$ws.Cells.Item(48,13).Value = $m47
$ws.Cells.Item(48,14).Value = $n47
$ws.Cells.Item(48,15).Value = $o47
$ws.Cells.Item(48,16).Value = $p47
$ws.Cells.Item(48,17).Value = $q47
$ws.Cells.Item(48,18).Value = $r47
